# Additional features, templates, etc
# The "VendorVLAN" config row is replaced with a "VendorTemplate" row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")
$ws.Activate()

$ws.Range("A7").Value = "VendorTemplate"
$ws.Range("C7").Value = "Template Name of where these devices should be sources (leave blank if none should be applied)"

# Leave the selection where the author left it when saving the file.
$ws.Range("B7").Select()
